$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 14929352
$ws.Range("J100").Value = 1431416.1
$ws.Range("L100").Value = 1431416.1
$ws.Range("N100").Value = -1432498.1
$ws.Range("H106").Value = 2052.2
$ws.Range("I106").Value = 1345.909
$ws.Range("J106").Value = 3994.5
$ws.Range("K106").Value = 1345.909
$ws.Range("L106").Value = 3994.5
$ws.Range("M106").Value = -714.9090000000001
$ws.Range("N106").Value = -5256.5
$ws.Range("H112").Value = 1043.7826
$ws.Range("I112").Value = 1080
$ws.Range("J112").Value = 1039.3658
$ws.Range("K112").Value = 3240
$ws.Range("L112").Value = 3118.0974
$ws.Range("M112").Value = -2132
$ws.Range("N112").Value = -5334.097400000001
$ws.Range("H129").Value = 1104.3549
$ws.Range("I129").Value = 465.7143
$ws.Range("J129").Value = 1290.625
$ws.Range("K129").Value = 1397.1429
$ws.Range("L129").Value = 3871.875
$ws.Range("M129").Value = 3602.8571
$ws.Range("N129").Value = -13871.875

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3179.6064
$ws.Range("I32").Value = 2069.5059
$ws.Range("K32").Value = 2069.5059
$ws.Range("M32").Value = -1782.5059
$ws.Range("H63").Value = 14769.9
$ws.Range("I63").Value = 14769.9
$ws.Range("K63").Value = 14769.9
$ws.Range("M63").Value = -14083.9
$ws.Range("H66").Value = 14769.9
$ws.Range("I66").Value = 14769.9
$ws.Range("K66").Value = 73849.5
$ws.Range("M66").Value = -70417.5
$ws.Range("H110").Value = 86963.86
$ws.Range("I110").Value = 200274
$ws.Range("K110").Value = 200274
$ws.Range("M110").Value = -198229
$ws.Range("H112").Value = 22837.4
$ws.Range("J112").Value = 22837.4
$ws.Range("L112").Value = 22837.4
$ws.Range("N112").Value = -25791.4
$ws.Range("H125").Value = 73315
$ws.Range("J125").Value = 73315
$ws.Range("L125").Value = 73315
$ws.Range("N125").Value = -83155
$ws.Range("H133").Value = 34628.332
$ws.Range("J133").Value = 34628.332
$ws.Range("L133").Value = 34628.332
$ws.Range("N133").Value = -39688.332

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H43").Value = 90000
$ws.Range("J43").Value = 90000
$ws.Range("L43").Value = 90000
$ws.Range("H64").Value = 418.82352
$ws.Range("J64").Value = 480.57144
$ws.Range("L64").Value = 480.57144
$ws.Range("N64").Value = -930.5714399999999
$ws.Range("H67").Value = 418.82352
$ws.Range("J67").Value = 480.57144
$ws.Range("L67").Value = 480.57144
$ws.Range("N67").Value = -2040.57144
$ws.Range("H112").Value = 63051.332
$ws.Range("J112").Value = 63051.332
$ws.Range("L112").Value = 63051.332
$ws.Range("N112").Value = -66005.33199999999
$ws.Range("N43").Value = -90362

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 844.2
$ws.Range("I10").Value = 305.25
$ws.Range("J10").Value = 3000
$ws.Range("K10").Value = 305.25
$ws.Range("L10").Value = 3000
$ws.Range("M10").Value = -166.25
$ws.Range("N10").Value = -3278

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 656927
$ws.Range("I68").Value = 1552437.8
$ws.Range("J68").Value = 1285.2322
$ws.Range("K68").Value = 4657313.4
$ws.Range("L68").Value = 3855.6966
$ws.Range("M68").Value = -4656502.4
$ws.Range("N68").Value = -5477.696599999999
$ws.Range("H71").Value = 656927
$ws.Range("I71").Value = 1552437.8
$ws.Range("J71").Value = 1285.2322
$ws.Range("K71").Value = 13971940.2
$ws.Range("L71").Value = 11567.0898
$ws.Range("M71").Value = -13967884.2
$ws.Range("N71").Value = -19679.0898
$ws.Range("H131").Value = 2981.9834
$ws.Range("J131").Value = 3525.9185
$ws.Range("L131").Value = 10577.7555
$ws.Range("N131").Value = -20657.7555
$ws.Range("H132").Value = 1450.0667
$ws.Range("I132").Value = 1622
$ws.Range("J132").Value = 1253.5714
$ws.Range("K132").Value = 14598
$ws.Range("L132").Value = 11282.1426
$ws.Range("M132").Value = -12068
$ws.Range("N132").Value = -16342.1426

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 14777.777
$ws.Range("I70").Value = 36333.332
$ws.Range("J70").Value = 4000
$ws.Range("K70").Value = 36333.332
$ws.Range("L70").Value = 4000
$ws.Range("M70").Value = -36063.332
$ws.Range("N70").Value = -4540
$ws.Range("H73").Value = 14777.777
$ws.Range("I73").Value = 36333.332
$ws.Range("J73").Value = 4000
$ws.Range("K73").Value = 36333.332
$ws.Range("L73").Value = 4000
$ws.Range("M73").Value = -35397.332
$ws.Range("N73").Value = -5872

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1254
$ws.Range("I46").Value = 1093.3334
$ws.Range("J46").Value = 1374.5
$ws.Range("K46").Value = 1093.3334
$ws.Range("L46").Value = 1374.5
$ws.Range("M46").Value = -905.3334
$ws.Range("N46").Value = -1750.5
$ws.Range("H130").Value = 40000
$ws.Range("J130").Value = 40000
$ws.Range("L130").Value = 40000
$ws.Range("N130").Value = -50040
$ws.Range("H136").Value = 5052126.5
$ws.Range("I136").Value = 1377.9584
$ws.Range("J136").Value = 18520790
$ws.Range("K136").Value = 4133.8752
$ws.Range("L136").Value = 55562370
$ws.Range("M136").Value = -1583.8752
$ws.Range("N136").Value = -55567470

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 70003.5
$ws.Range("J18").Value = 70003.5
$ws.Range("L18").Value = 70003.5
$ws.Range("N18").Value = -70349.5
$ws.Range("H122").Value = 69045.2
$ws.Range("I122").Value = 145386.86
$ws.Range("J122").Value = 2246.25
$ws.Range("K122").Value = 436160.58
$ws.Range("L122").Value = 6738.75
$ws.Range("M122").Value = -433710.58
$ws.Range("N122").Value = -11638.75
$ws.Range("H123").Value = 43542
$ws.Range("J123").Value = 43542
$ws.Range("L123").Value = 43542
$ws.Range("N123").Value = -53342
$ws.Range("H124").Value = 19250
$ws.Range("J124").Value = 19250
$ws.Range("L124").Value = 19250
$ws.Range("H132").Value = 2014.4822
$ws.Range("I132").Value = 1369.4333
$ws.Range("J132").Value = 2758.7693
$ws.Range("K132").Value = 4108.2999
$ws.Range("L132").Value = 8276.3079
$ws.Range("M132").Value = -1578.2999
$ws.Range("N132").Value = -13336.3079
$ws.Range("N124").Value = -29070
